$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Image front" / "Image back" columns (both the table header cells
# and, since this is a structured Table, the bound ListObject column names)
# to "ImageFront" / "ImageBack".
$ws.Range("E1").Value = "ImageFront"
$ws.Range("F1").Value = "ImageBack"

# Row 6 (Board #3) had its front/back image filenames swapped.
$ws.Range("E6").Value = "board3-back.jpg"
$ws.Range("F6").Value = "board3-front.jpg"

# Update the saved selection to match the authored state.
$ws.Range("E13").Select() | Out-Null
